# Data processing and sub-categorization
#
# Re-layout the Gas "Type" / "Form" table:
#   - Hydrogen keeps a single "Not Labelled" sub-row (row 2)
#   - A new "Acetylene" type is introduced with 6 Form sub-rows
#     (rows 3-8, merged A3:A8)
#   - "Other (Specify)" becomes its own Type group with 2 Form sub-rows
#     (rows 9-10, merged A9:A10, new rows appended at the bottom)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Break apart the existing "Type" merges -- writing to a non-anchor cell
#    of a merged range is a no-op in Excel, so the old A2:A6 / A7:A8 merges
#    must come apart before the cells underneath get new content.
$ws.Range("A2:A6").UnMerge()
$ws.Range("A7:A8").UnMerge()

# 2) Clear formatting on the rows that will be re-merged so that the
#    upcoming Merge() calls don't try to split the existing thin border
#    across the newly-merged block (format gets reapplied via PasteSpecial
#    in step 5).
$ws.Range("A3:A10").ClearFormats()

# 3) Re-merge the "Type" column into the new groupings.
$ws.Range("A3:A8").Merge()
$ws.Range("A9:A10").Merge()

# 4) Write the new cell values, row by row.
# Row 2: Hydrogen / Not Labelled
$ws.Range("A2").Value = "Hydrogen"
$ws.Range("B2").Value = "Not Labelled"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.2

# Row 3: Acetylene / Jet from a pressurized source
$ws.Range("A3").Value = "Acetylene"
$ws.Range("B3").Value = "Jet from a pressurized source"
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 1.7

# Row 4: Pressurized in a container
$ws.Range("B4").Value = "Pressurized in a container"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 0.8

# Row 5: Other (Specify)
$ws.Range("B5").Value = "Other (Specify)"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 0.6

# Row 6: Ambient pressure, within a compartment
$ws.Range("B6").Value = "Ambient pressure, within a compartment"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0.4

# Row 7: Ambient pressure, within a component
$ws.Range("B7").Value = "Ambient pressure, within a component"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.4

# Row 8: Not Labelled (closes the Acetylene merged block)
$ws.Range("B8").Value = "Not Labelled"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.2

# Row 9: Other (Specify) / Other (Specify) -- new row
$ws.Range("A9").Value = "Other (Specify)"
$ws.Range("B9").Value = "Other (Specify)"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0.4

# Row 10: Pressurized in a container -- new row
$ws.Range("B10").Value = "Pressurized in a container"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0.2

# 5) Re-apply the bold/centered/bordered "Type"+"Form" header style to the
#    A and B columns for the rows that were touched above, matching the
#    formatting already used by the rest of the table (copied from A1/B1).
$ws.Range("A1").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B2:B10").PasteSpecial(-4122)

$excel.CutCopyMode = $false
